# Adjusted risk calc formula
# Updates several cells in the Nahmint watershed risk table (Table 1):
#   - LF70 row (table row 13): Total Risk 4 -> 6, Current Risk L -> M
#   - LF2  row (table row 14): Total Risk 2 -> 4, Future Risk  VL -> L
#   - LF5  row (table row 15): Total Risk 2 -> 4, Future Risk  VL -> L
#   - LF8  row (table row 16): Rank 13 -> 15
#   - LF21 row (table row 17): Rank 13 -> 15
#   - LF38 row (table row 18): Rank 13 -> 15
#   - LF39 row (table row 19): Rank 13 -> 15
#
# Cell(row, col) is used instead of Find/Replace because Find.Execute in this
# runtime operates on the whole document regardless of the range/cell it is
# invoked on, which would incorrectly replace matching text elsewhere in the
# document (e.g. every "L" or digit). Assigning directly to Range.Text keeps
# the edit scoped to the single target cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Columns: 1=Watershed, 2=LF, 3=Rank, 4=Total Risk, 5=Current Risk, 6=Future Risk

# Row 13 = LF70
$t.Cell(13, 4).Range.Text = "6"
$t.Cell(13, 5).Range.Text = "M"

# Row 14 = LF2
$t.Cell(14, 4).Range.Text = "4"
$t.Cell(14, 6).Range.Text = "L"

# Row 15 = LF5
$t.Cell(15, 4).Range.Text = "4"
$t.Cell(15, 6).Range.Text = "L"

# Row 16 = LF8
$t.Cell(16, 3).Range.Text = "15"

# Row 17 = LF21
$t.Cell(17, 3).Range.Text = "15"

# Row 18 = LF38
$t.Cell(18, 3).Range.Text = "15"

# Row 19 = LF39
$t.Cell(19, 3).Range.Text = "15"
